$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 6 ("grandes regiões e unidades da federação" header row),
# which shifts all subsequent rows (and their numeric data) up by one,
# and removes the now-unused shared string from the workbook.
$ws.Rows.Item(6).Delete()
